# Append-refresh of the "ランサーズ" (案件情報) sheet: 2025-10-24 12:40 JST run.
# The scraper re-ran, re-sorted by priority score (desc) and produced a new
# top-12 list; every "取得日時" timestamp moves to the new run time, several
# rows keep the same listing (shifted position), a few listings are new, and
# one listing (Google Play Console clouse-test, score 10) aged out of the
# top results entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-10-24 12:40:12"

# row, title, category, price, deadline, url, score, skills(or $null)
$rowsData = @(
    @(2,  "自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5408668", 305, "🔥Python ◆開発 ○PHP"),
    @(3,  "製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419380", 298, "🔥AI,Ai"),
    @(4,  "【急募】経験豊富な業務システム開発パートナーを募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419860", 125, "◆開発,システム開発"),
    @(5,  "【低コスト】住宅リフォーム見積依頼自動化システム構築", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5412955", 110, "◆自動化"),
    @(6,  "Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419636", 85, "★Java"),
    @(7,  "UIPATHのシステムの開発", "システム開発", "10,000 円 ~", "期限情報なし", "https://www.lancers.jp/work/detail/5419904", 75, "◆開発"),
    @(8,  "IB報酬を得るための高性能EA開発依頼", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419587", 68, "◆開発"),
    @(9,  "クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419638", 38, "◇管理"),
    @(10, "【緊急】ロリポップ Wordpress リダイレクトハッキング復旧依頼", "システム開発", "20,000 円 ~ 30,000 円 / 募集期間 5 日、取引期間 0 日", "期限情報なし", "https://www.lancers.jp/work/detail/5419656", 25, "○WordPress"),
    @(11, "【Braze経験者募集】CRM/マーケティングオートメーション支援(中級者以上)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419838", 25, $null),
    @(12, "【KARTE経験者募集】CX改善/Web接客施策の設計・実装(中級者以上)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419829", 25, $null)
)

# Drop every pre-existing hyperlink relationship up front -- row positions
# are being rewritten wholesale below, and this engine does not re-target
# hyperlink refs when cells move, so the clean way is: clear, then re-add
# in the correct final order.
$ws.Hyperlinks.Delete()

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
    $ws.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
}

# Column D ("価格") widened 28 -> 41 characters to fit the longer strings
# (e.g. the new "募集期間 5 日、取引期間 0 日" wording).
$ws.Columns.Item(4).ColumnWidth = 40.17
